# "attempted one leetcode qn" - add the Maximum SubArray (Kadane's algorithm) row
# to the LeetCode tracker sheet, then tidy up row heights / column widths /
# the pre-formatted blank rows below the table the same way the author's
# Excel session ended up saving them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row -----------------------------------------------------
# Row 4 already holds the "Merge Two Sorted Lists" entry; copy its
# formatting down into row 5 for the new question before filling it in.
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A5:C5").PasteSpecial(-4122) | Out-Null

$ws.Range("A5").Value = "Maximum SubArray"
$ws.Range("B5").Value = "Easy"
$ws.Range("C5").Value = "Kadane's algorithm. Loop through once, if the currentSum < 0 initialize back to 0. need to self check cause sometimes it may not be less than 0 but can be smaller than previous combination"

# --- row heights --------------------------------------------------------
# Rows 2-4 (existing questions) shrink slightly to the new uniform height,
# and the new row 5 plus the already-present blank rows 6-14 take it too.
$ws.Range("A2:C16").RowHeight = 25.05

# A further block of blank rows (15 through 50) was pre-formatted with the
# same row height, ready for future entries.
$ws.Range("A17:A50").EntireRow.RowHeight = 25.05

# --- column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.8333333333333
$ws.Columns.Item(2).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 153.333333333333

# --- selection -------------------------------------------------------
[void]$ws.Range("C4").Select()
